# Adapt column header formatting to respective input file names.
# "_old"/"_new" suffixes become "_FV2310"/"_FV2404" respectively, the
# header row A1:U1 is turned into a proper Excel Table (with AutoFilter),
# and the top header row is frozen so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header cells (A1:U1) -----------------------------------
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2) Turn the header + data range into a real Excel Table --------------
$dataRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)   # xlSrcRange, hdr=xlYes
$tbl.Name = "Table1"

# --- 3) Freeze the header row (split below row 1) --------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
